function Get-CellText($cell) {
    # Cell text always ends with a cell-mark / paragraph-mark control
    # character; strip those so plain string comparisons work.
    return $cell.Range.Text.TrimEnd([char]7, [char]13)
}

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "grade_id" foreign-key row into the "competition"
#    data-dictionary table, directly above the "title" row.
# ------------------------------------------------------------------
$competitionTable = $null
foreach ($tbl in $d.Tables) {
    if ((Get-CellText $tbl.Cell(1, 1)) -eq "competition") {
        $competitionTable = $tbl
        break
    }
}

$titleRow = $null
foreach ($r in $competitionTable.Rows) {
    if ((Get-CellText $r.Cells.Item(2)) -eq "title") {
        $titleRow = $r
        break
    }
}

$newRow = $competitionTable.Rows.Add($titleRow)  # inserts a blank row just above "title"
$newRow.Height = 15.75                           # 315 twips, matching its neighbouring rows

$newRow.Cells.Item(1).Range.Text = "PK,FK1"
$newRow.Cells.Item(2).Range.Text = "grade_id"
$newRow.Cells.Item(3).Range.Text = "Integer"
# 4th cell (Description) is left blank, as in the source edit.

# ------------------------------------------------------------------
# 2. Re-apply the "field" table-title text so the stale cached
#    lastRenderedPageBreak marker in front of it is cleared out.
# ------------------------------------------------------------------
$fieldTable = $null
foreach ($tbl in $d.Tables) {
    if ((Get-CellText $tbl.Cell(1, 1)) -eq "field") {
        $fieldTable = $tbl
        break
    }
}
$fieldTable.Cell(1, 1).Range.Text = "field"
